$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet already contains two duplicated "a1"/"b2" records in rows 7-8
# and 9-10. This change appends another duplicate pair as rows 11-12,
# extending the used range from A1:H10 to A1:H12.
#
# Copy rows 9 and 10 (source of the duplicate pattern) down to the new
# rows 11 and 12, using the single-call Copy(Destination) form so that
# values, number formatting and the bordered/bold/centered cell style
# (applied to column A) are carried over exactly as in the original rows.
$ws.Range("A9:H9").Copy($ws.Range("A11:H11"))
$ws.Range("A10:H10").Copy($ws.Range("A12:H12"))
